$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data rows 2-11: employee_id, employee_name, department, absence_reason, absence_duration, absence_date, salary
$data = @(
    @(90206, "Ravi Lucca Rezende", "TI", "Consulta medica", 2, 45091, 3855.91),
    @(13657, "Calebe Campos", "Recursos Humanos", "Outros", 2, 45104, 7163.88),
    @(69514, "Anna Liz Azevedo", "Financeiro", "Problemas pessoais", 4, 45095, 3385.71),
    @(27382, "Otávio Moreira", "TI", "Consulta medica", 7, 45080, 4971.99),
    @(36251, "Sra. Emanuella Montenegro", "Operacoes", "Consulta medica", 1, 45095, 4142.27),
    @(83233, "Dra. Mariana da Cruz", "Vendas", "Doenca", 6, 45095, 8069.8),
    @(49594, "Joana Camargo", "Financeiro", "Outros", 6, 45106, 3438.49),
    @(51703, "Isabella Pereira", "Engenharia", "Problemas pessoais", 4, 45089, 6274.75),
    @(20352, "Jade Moraes", "Juridico", "Consulta medica", 7, 45095, 5395.39),
    @(34462, "Manuela da Cunha", "TI", "Consulta medica", 6, 45089, 3204.64)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowIndex = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $ws.Cells.Item($rowIndex, 7).Value = $row[6]
}
